# use K+ for combined dropped third strike
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("codes")

# Row 31: Droppped Third Strike -> add F31 = "K+" (new VAR3 value)
$ws.Range("F31").Value = "K+"
$ws.Range("F31").HorizontalAlignment = $ws.Range("E31").HorizontalAlignment

# Row 32: Passed Ball (on Dropped Third Strike)
# VAR1 (D32) now matches the Code (PB) instead of "Kd"; VAR3 (F32) becomes "K+"
$ws.Range("D32").Value = "PB"
$ws.Range("F32").Value = "K+"

# Row 33: Wild Pitch (on Dropped Third Strike)
# VAR1 (D33) now matches the Code (WP) instead of "Kd"; VAR3 (F33) becomes "K+"
$ws.Range("D33").Value = "WP"
$ws.Range("F33").Value = "K+"

# Row 43: Passed Ball - VAR1 (D43) now matches Code (PB) instead of "E";
# the separate VAR3 (F43) cell is removed entirely
$ws.Range("D43").Value = "PB"
$ws.Range("F43").Clear()

# Row 46: Wild Pitch - VAR1 (D46) now matches Code (WP) instead of "E";
# the separate VAR3 (F46) cell is removed entirely
$ws.Range("D46").Value = "WP"
$ws.Range("F46").Clear()

# Restore the view state (frozen pane top-left cell and active selection)
$ws.Application.ActiveWindow.Panes.Item(4).ScrollRow = 16
$ws.Range("E28").Select()
